$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2342.8572
$ws.Range("J40").Value = 3466.6667
$ws.Range("L40").Value = 3466.6667
$ws.Range("N40").Value = -3816.6667
$ws.Range("H47").Value = 12500
$ws.Range("I47").Value = 12500
$ws.Range("J47").Value = 0
$ws.Range("K47").Value = 12500
$ws.Range("L47").Value = 0
$ws.Range("M47").Value = -11528
$ws.Range("N47").ClearContents()
$ws.Range("H92").Value = 1570.6923
$ws.Range("I92").Value = 1519.0435
$ws.Range("J92").Value = 1966.6666
$ws.Range("K92").Value = 1519.0435
$ws.Range("L92").Value = 1966.6666
$ws.Range("M92").Value = -271.0435
$ws.Range("N92").Value = -4462.6666
$ws.Range("H106").Value = 3220
$ws.Range("I106").Value = 3025
$ws.Range("K106").Value = 3025
$ws.Range("M106").Value = -2394
$ws.Range("H116").Value = 4662.375
$ws.Range("I116").Value = 5259.8
$ws.Range("J116").Value = 3666.6667
$ws.Range("K116").Value = 5259.8
$ws.Range("L116").Value = 3666.6667
$ws.Range("M116").Value = -1817.8
$ws.Range("N116").Value = -10550.6667
$ws.Range("H135").Value = 5600.2
$ws.Range("I135").Value = 4000.3333
$ws.Range("J135").Value = 8000
$ws.Range("K135").Value = 36002.9997
$ws.Range("L135").Value = 72000
$ws.Range("M135").Value = -33467.9997
$ws.Range("N135").Value = -77070
$ws.Range("H137").Value = 1364
$ws.Range("I137").Value = 986.2857
$ws.Range("J137").Value = 2025
$ws.Range("K137").Value = 2958.8571
$ws.Range("L137").Value = 6075
$ws.Range("M137").Value = -408.8571000000002
$ws.Range("N137").Value = -11175

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4947.9
$ws.Range("I32").Value = 4201.7407
$ws.Range("K32").Value = 4201.7407
$ws.Range("M32").Value = -3914.7407
$ws.Range("H88").Value = 3625.7273
$ws.Range("I88").Value = 3261.2
$ws.Range("J88").Value = 3929.5
$ws.Range("K88").Value = 3261.2
$ws.Range("L88").Value = 3929.5
$ws.Range("M88").Value = -2855.2
$ws.Range("N88").Value = -4741.5
$ws.Range("H91").Value = 3625.7273
$ws.Range("I91").Value = 3261.2
$ws.Range("J91").Value = 3929.5
$ws.Range("K91").Value = 3261.2
$ws.Range("L91").Value = 3929.5
$ws.Range("M91").Value = -1857.2
$ws.Range("N91").Value = -6737.5
$ws.Range("H97").Value = 743.65216
$ws.Range("I97").Value = 721.2105
$ws.Range("J97").Value = 850.25
$ws.Range("K97").Value = 721.2105
$ws.Range("L97").Value = 850.25
$ws.Range("M97").Value = -225.2105
$ws.Range("N97").Value = -1842.25
$ws.Range("H101").Value = 20750
$ws.Range("J101").Value = 20750
$ws.Range("L101").Value = 20750
$ws.Range("N101").Value = -27240

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H55").Value = 54945
$ws.Range("J55").Value = 54945
$ws.Range("L55").Value = 54945
$ws.Range("N55").Value = -55491
$ws.Range("H94").Value = 794.9091
$ws.Range("I94").Value = 810.5333000000001
$ws.Range("J94").Value = 761.4286
$ws.Range("K94").Value = 810.5333000000001
$ws.Range("L94").Value = 761.4286
$ws.Range("M94").Value = -359.5333000000001
$ws.Range("N94").Value = -1663.4286
$ws.Range("H100").Value = 13750
$ws.Range("J100").Value = 13750
$ws.Range("L100").Value = 13750
$ws.Range("N100").Value = -15914
$ws.Range("H105").Value = 6785.8823
$ws.Range("I105").Value = 6677.5
$ws.Range("J105").Value = 7046
$ws.Range("K105").Value = 6677.5
$ws.Range("L105").Value = 7046
$ws.Range("M105").Value = -4930.5
$ws.Range("N105").Value = -10540
$ws.Range("H107").Value = 18884.875
$ws.Range("I107").Value = 2099.6667
$ws.Range("J107").Value = 28956
$ws.Range("K107").Value = 2099.6667
$ws.Range("L107").Value = 28956
$ws.Range("M107").Value = -179.6667000000002
$ws.Range("N107").Value = -32796
$ws.Range("H114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("L114").Value = 0
$ws.Range("N114").ClearContents()
$ws.Range("H134").Value = 93186.82000000001
$ws.Range("I134").Value = 2338.6667
$ws.Range("J134").Value = 502003.5
$ws.Range("K134").Value = 7016.000100000001
$ws.Range("L134").Value = 1506010.5
$ws.Range("M134").Value = -4481.000100000001
$ws.Range("N134").Value = -1511080.5

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4088.7917
$ws.Range("I31").Value = 4310.85
$ws.Range("J31").Value = 2978.5
$ws.Range("K31").Value = 4310.85
$ws.Range("L31").Value = 2978.5
$ws.Range("M31").Value = -4015.85
$ws.Range("N31").Value = -3568.5
$ws.Range("H34").Value = 4088.7917
$ws.Range("I34").Value = 4310.85
$ws.Range("J34").Value = 2978.5
$ws.Range("K34").Value = 4310.85
$ws.Range("L34").Value = 2978.5
$ws.Range("M34").Value = -4108.85
$ws.Range("N34").Value = -3382.5
$ws.Range("H107").Value = 363.97144
$ws.Range("I107").Value = 371
$ws.Range("J107").Value = 361.53845
$ws.Range("K107").Value = 371
$ws.Range("L107").Value = 361.53845
$ws.Range("M107").Value = 1549
$ws.Range("N107").Value = -4201.53845
$ws.Range("H122").Value = 767.2
$ws.Range("I122").Value = 709
$ws.Range("J122").Value = 1000
$ws.Range("K122").Value = 2127
$ws.Range("L122").Value = 3000
$ws.Range("M122").Value = 323
$ws.Range("N122").Value = -7900
$ws.Range("H132").Value = 1455.9117
$ws.Range("I132").Value = 896.65216
$ws.Range("J132").Value = 2625.2727
$ws.Range("K132").Value = 2689.95648
$ws.Range("L132").Value = 7875.8181
$ws.Range("M132").Value = -159.9564799999998
$ws.Range("N132").Value = -12935.8181

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H121").Value = 945.4286
$ws.Range("I121").Value = 968
$ws.Range("J121").Value = 938.375
$ws.Range("K121").Value = 2904
$ws.Range("L121").Value = 2815.125
$ws.Range("M121").Value = -1594
$ws.Range("N121").Value = -5435.125
$ws.Range("H130").Value = 4065.7144
$ws.Range("I130").Value = 1833.3334
$ws.Range("J130").Value = 5740
$ws.Range("K130").Value = 5500.0002
$ws.Range("L130").Value = 17220
$ws.Range("M130").Value = -480.0002000000004
$ws.Range("N130").Value = -27260
$ws.Range("H131").Value = 783.74
$ws.Range("I131").Value = 265.36365
$ws.Range("J131").Value = 847.80896
$ws.Range("K131").Value = 796.09095
$ws.Range("L131").Value = 2543.42688
$ws.Range("M131").Value = 4243.90905
$ws.Range("N131").Value = -12623.42688

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H32").Value = 20290
$ws.Range("J32").Value = 20290
$ws.Range("L32").Value = 20290
$ws.Range("N32").Value = -20882
$ws.Range("H42").Value = 25017144
$ws.Range("J42").Value = 25017144
$ws.Range("L42").Value = 25017144
$ws.Range("N42").Value = -25018114
$ws.Range("H45").Value = 20326
$ws.Range("J45").Value = 20326
$ws.Range("L45").Value = 20326
$ws.Range("N45").Value = -21444
$ws.Range("H51").Value = 21333.334
$ws.Range("I51").Value = 15000
$ws.Range("J51").Value = 24500
$ws.Range("K51").Value = 15000
$ws.Range("L51").Value = 24500
$ws.Range("M51").Value = -14491
$ws.Range("N51").Value = -25518
$ws.Range("H97").Value = 1647.2222
$ws.Range("I97").Value = 1754.2142
$ws.Range("J97").Value = 1272.75
$ws.Range("K97").Value = 1754.2142
$ws.Range("L97").Value = 1272.75
$ws.Range("M97").Value = -1258.2142
$ws.Range("N97").Value = -2264.75
$ws.Range("H115").Value = 25017144
$ws.Range("J115").Value = 25017144
$ws.Range("L115").Value = 25017144
$ws.Range("N115").Value = -25019494
$ws.Range("H126").Value = 3103.25
$ws.Range("I126").Value = 4970.6665
$ws.Range("J126").Value = 1982.8
$ws.Range("K126").Value = 14911.9995
$ws.Range("L126").Value = 5948.4
$ws.Range("M126").Value = -12441.9995
$ws.Range("N126").Value = -10888.4

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1112.6364
$ws.Range("I46").Value = 2397.5
$ws.Range("J46").Value = 378.42856
$ws.Range("K46").Value = 2397.5
$ws.Range("L46").Value = 378.42856
$ws.Range("M46").Value = -2209.5
$ws.Range("N46").Value = -754.4285600000001
$ws.Range("H68").Value = 2770.2
$ws.Range("I68").Value = 3200.2856
$ws.Range("J68").Value = 1766.6666
$ws.Range("K68").Value = 3200.2856
$ws.Range("L68").Value = 1766.6666
$ws.Range("M68").Value = -2451.2856
$ws.Range("N68").Value = -3264.6666
$ws.Range("H71").Value = 2770.2
$ws.Range("I71").Value = 3200.2856
$ws.Range("J71").Value = 1766.6666
$ws.Range("K71").Value = 16001.428
$ws.Range("L71").Value = 8833.333000000001
$ws.Range("M71").Value = -12257.428
$ws.Range("N71").Value = -16321.333

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H43").Value = 37500
$ws.Range("J43").Value = 38333.332
$ws.Range("L43").Value = 38333.332
$ws.Range("N43").Value = -38631.332
$ws.Range("H49").Value = 9000
$ws.Range("J49").Value = 9000
$ws.Range("L49").Value = 9000
$ws.Range("N49").Value = -9460
$ws.Range("H109").Value = 19811.111
$ws.Range("J109").Value = 19811.111
$ws.Range("L109").Value = 19811.111
